$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Repeticion indefinida, acumuladores," + _GoBack bookmark + " contadores"
# becomes a single run "Repeticion indefinida, acumuladores, contadores" with
# the leftover _GoBack bookmark removed. Word's Find treats bookmarks as
# zero-width, so searching across the old run/bookmark/run boundary and
# replacing with the same (now contiguous) text coalesces everything into one
# run and drops the now-empty bookmark.
# ---------------------------------------------------------------------------
$found1 = $d.Content.Find.Execute(", acumuladores, contadores", $true, $false, $false, $false, $false, $true, 1, $false, ", acumuladores, contadores", 2)

# ---------------------------------------------------------------------------
# Edit 2: "\u00a1Acertaste! Te ha costado  3 intentos" (double space) becomes
# two runs - "\u00a1Acertaste! Te ha costado " and "3 intentos" - split by a
# fresh _GoBack bookmark, and the double space collapses to a single space.
# ---------------------------------------------------------------------------
$found2 = $d.Content.Find.Execute("Te ha costado  3 intentos", $true, $false, $false, $false, $false, $true, 1, $false, "Te ha costado 3 intentos", 2)

$rng2 = $d.Content
$found2b = $rng2.Find.Execute("3 intentos")
$rng2.Collapse(1)
$d.Bookmarks.Add("_GoBack", $rng2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 3: " y la cantidad de numeros ingreso sin considerar el -1" becomes
# three runs - " y la cantidad de numeros ", "que " and
# "ingreso sin considerar el -1" - by inserting the word "que " in the middle.
# Temporary bookmarks are dropped at both sides of the insertion point (and
# the pre-existing run boundary to its left) to stop the engine's
# same-formatting run merge from swallowing the split; they are deleted again
# once the text is in place, leaving only the three plain runs behind.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute(" y la cantidad de números ingresó sin considerar el -1")
$origStart = $rng3.Start

$leftBoundary = $origStart
$splitPos = $origStart + 26   # length of " y la cantidad de números "

$wallLeft = $d.Bookmarks.Add("zzwallLeft", $d.Range($leftBoundary, $leftBoundary))
$wallMid = $d.Bookmarks.Add("zzwallMid", $d.Range($splitPos, $splitPos))

$insertRng = $d.Range($splitPos, $splitPos)
$insertRng.InsertBefore("que ")

$wallPreQue = $d.Bookmarks.Add("zzwallPreQue", $d.Range($splitPos, $splitPos))

$d.Bookmarks("zzwallLeft").Delete()
$d.Bookmarks("zzwallMid").Delete()
$d.Bookmarks("zzwallPreQue").Delete()

Write-Output "edit1=$found1 edit2=$found2/$found2b edit3=$found3"
